# Automatic update of files.
# The underlying data rows got re-synced from source: each row's content
# (species/taxon identity, coordinates, etc.) was rotated among the rows
# 2-9 (matched by the "Id" / A column), GPS easting/northing (Q/R) were
# rounded to whole metres, and the Starttid/Sluttid (Z/AB) columns were
# dropped.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-RowData {
    param($Row, $A, $B, $D, $E, $F, $G, $H, $I, $M, $Q, $R)

    $ws.Range("A$Row").Value = $A
    $ws.Range("B$Row").Value = $B
    if ($D -ne $null) { $ws.Range("D$Row").Value = $D }
    $ws.Range("E$Row").Value = $E
    $ws.Range("F$Row").Value = $F
    $ws.Range("G$Row").Value = $G
    $ws.Range("H$Row").Value = $H

    if ($I -ne $null) {
        # Force text storage so a numeric-looking value like "3" keeps
        # matching the source file's inlineStr cell type instead of being
        # auto-coerced to a number by the normal Excel input-parsing rules.
        $ws.Range("I$Row").NumberFormat = "@"
        $ws.Range("I$Row").Value = $I
    } else {
        $ws.Range("I$Row").ClearContents()
    }

    if ($M -ne $null) {
        $ws.Range("M$Row").Value = $M
    } else {
        $ws.Range("M$Row").ClearContents()
    }

    $ws.Range("Q$Row").Value = $Q
    $ws.Range("R$Row").Value = $R

    # Starttid / Sluttid no longer populated.
    $ws.Range("Z$Row").ClearContents()
    $ws.Range("AB$Row").ClearContents()
}

# Row 2 <- old row 7 (111739317 / Skrovellav)
Set-RowData 2 111739317 78579 "NT" 2081 "Skrovellav" "Lobaria scrobiculata" "(Scop.) DC." $null $null 573912 7172648

# Row 3 <- old row 6 (111739307 / Talltita)
Set-RowData 3 111739307 56543 "NT" 103021 "Talltita" "Poecile montanus" "(Conrad von Baldenstein, 1827)" "3" "födosökande" 573961 7172501

# Row 4 <- old row 9 (111739313 / Rödbrun blekspik)
Set-RowData 4 111739313 73701 "NT" 1467 "Rödbrun blekspik" "Sclerophora coniophaea" "(Norman) J.Mattsson & Middelb." $null $null 574025 7172443

# Row 5 <- old row 3 (111739311 / Garnlav)
Set-RowData 5 111739311 77515 "NT" 6425 "Garnlav" "Alectoria sarmentosa" "(Ach.) Ach." $null $null 574012 7172473

# Row 6 <- old row 5 (111739315 / Stuplav)
Set-RowData 6 111739315 78605 "LC" 6462 "Stuplav" "Nephroma bellum" "(Spreng.) Tuck." $null $null 573905 7172637

# Row 7 <- old row 4 (111739306 / Tretåig hackspett)
Set-RowData 7 111739306 56398 $null 100109 "Tretåig hackspett" "Picoides tridactylus" "(Linnaeus, 1758)" $null "äldre spår" 573906 7172521

# Row 8 <- old row 8 (111739316 / Lunglav) - identity unchanged, only Q/R rounding + Z/AB removed
$ws.Range("Q8").Value = 573905
$ws.Range("R8").Value = 7172637
$ws.Range("Z8").ClearContents()
$ws.Range("AB8").ClearContents()

# Row 9 <- old row 2 (111739309 / Korallblylav)
Set-RowData 9 111739309 78536 "LC" 229497 "Korallblylav" "Parmeliella triptophylla" "(Ach.) Müll.Arg." $null $null 574011 7172434

Write-Host "Row rotation applied."
